# Add a "Span (km)" data column to the QSM output sheet.
# The new column is inserted as the 4th column (between "Q (dB)" and
# "Compensation (%)"), pushing the existing D:H columns to E:I, and every
# data row gets a Span value of 30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D - shifts old D:H (Compensation %, fiberAeff_1/2,
# fiberAlphadB_1/2) one column to the right, to E:I.
$ws.Columns("D:D").Insert()

# Header for the newly inserted column.
$ws.Range("D1").Value = "Span (km)"

# Fill the new column's data rows (2-14) with the span value.
$ws.Range("D2:D14").Value = 30
